# "new update day 08/10/2024"
# The sheet previously held a "link" header in A1 followed by 30 Facebook
# post URLs in A2:A31. This update clears out all of the scraped URL rows,
# leaving just the header, and widens the workbook window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 30 data rows (A2:A31) entirely - this also prunes the now
# unused shared strings (the URLs) so only "link" remains in sharedStrings.
$ws.Range("A2:A31").EntireRow.Delete()

# Move / restore the active selection to C9, matching the saved view state.
$ws.Range("C9").Select()

# Widen the saved window.
$excel.ActiveWindow.Width = 27945
